$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "296.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.73%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.08%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.042"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.61%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07567"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.72%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.603"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.03%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9283"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.45%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.402"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.87%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1217"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.78%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1834"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "6.52%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08947"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.08%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04017"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.84%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1051"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.12%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001280"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.18%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005866"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.17%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.361"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.88%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.402"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.86%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.11%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.905"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.58%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.05%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3004"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.21%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04059"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.50%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003958"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "4.09%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.08%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02406"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4.19%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05206"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.85%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006041"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-7.61%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007780"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.27%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1331"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.34%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007568"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.66%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007860"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "11.14%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.2970"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.72%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006790"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.72%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.20%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "165.66%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004207"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.13%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.20%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.20%"
